$d = $word.ActiveDocument

$pairs = @(
    @("196×9=", "432×7="),
    @("454×8=", "204×3="),
    @("894×4=", "780×3="),
    @("914×2=", "606×3="),
    @("953×2=", "702×3="),
    @("880×5=", "116×7="),
    @("356×3=", "495×9="),
    @("645×9=", "734×7="),
    @("407×3=", "977×5="),
    @("369×5=", "678×5="),
    @("441×9=", "755×5="),
    @("482×9=", "851×5="),
    @("975×4=", "259×7="),
    @("611×8=", "911×7="),
    @("297×4=", "547×3="),
    @("569×5=", "313×4="),
    @("942×2=", "611×4="),
    @("866×4=", "221×8="),
    @("125×3=", "110×3="),
    @("226×8=", "384×9="),
    @("341×7=", "471×8="),
    @("421×5=", "612×8="),
    @("734×5=", "463×8="),
    @("478×9=", "362×9="),
    @("445×4=", "893×4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
